# Natmi following Dr Hou advice
# Replace the two-row LR-pair result table with the updated 4-row table
# (sending/target cluster pairs now include "ECs" in addition to "sCs",
# and every numeric column has been recomputed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs","Tg","Asgr1","FAPs",1,0.3333333333333333,1.030842333333333,3.092527,0.7593378091926586,0.7593378091926584,3,1,1.164453,3.493359,0.7897445044585573,0.7897445044585573,1.200367447577,10.803307028193,0.5996828618375026,0.5996828618375025),
    @("ECs","Tg","Asgr1","sCs",1,0.3333333333333333,1.030842333333333,3.092527,0.7593378091926586,0.7593378091926584,3,1,0.310015,0.930045,0.2102554955414427,0.2102554955414428,0.3195765859683333,2.876189273715,0.1596549473551559,0.1596549473551559),
    @("sCs","Tg","Asgr1","FAPs",3,1,0.326712,0.980136,0.2406621908073416,0.2406621908073415,3,1,1.164453,3.493359,0.7897445044585573,0.7897445044585573,0.380440768536,3.423966916824,0.1900616426210547,0.1900616426210547),
    @("sCs","Tg","Asgr1","sCs",3,1,0.326712,0.980136,0.2406621908073416,0.2406621908073415,3,1,0.310015,0.930045,0.2102554955414427,0.2102554955414428,0.10128562068,0.91157058612,0.05060054818628684,0.05060054818628684)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowIndex = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $ws.Cells.Item($rowIndex, $j + 1).Value = $rowValues[$j]
    }
}
